# "Started adding all tables" - rename several worksheets to the final
# singular/plural naming convention used across the workbook, and start
# fleshing out the Users table with the columns the app actually needs
# (FavouriteTeamID, IsAuthenticated, IsActive, IsAdmin).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename worksheets
# ---------------------------------------------------------------------
$wb.Worksheets.Item("PredictionPeriodInfo").Name   = "PredictionPeriods"
$wb.Worksheets.Item("PredictionsByWeek").Name      = "WeekPredictions"
$wb.Worksheets.Item("PrivateLeagueInfo").Name      = "PrivateLeagues"
$wb.Worksheets.Item("ScoringSettingsInfo").Name    = "ScoringSettings"
$wb.Worksheets.Item("PrivateLeagueMembership").Name = "PrivateLeagueMembers"

# ---------------------------------------------------------------------
# 2. Users sheet: insert a new "IsActive" column and rename a few others
#    so the table reads: UserID, Username, EmailAddress, FavouriteTeamID,
#    HashedPassword, IsAuthenticated, IsActive, IsAdmin, UserCreated
# ---------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

# Make room for the new "IsActive" column between "Authenticated" (F) and
# "Admin" (G, about to become "IsAdmin").
$users.Columns("G").Insert()

$users.Range("D1").Value = "FavouriteTeamID"
$users.Range("D2").Value = "NULL"

$users.Range("F1").Value = "IsAuthenticated"

$users.Range("G1").Value = "IsActive"
$users.Range("G2").Value = 1

$users.Range("H1").Value = "IsAdmin"

$users.Range("G3").Select()

# ---------------------------------------------------------------------
# 3. Leave off on the PeriodPredictions sheet, ready to add the next row
# ---------------------------------------------------------------------
$periodPredictions = $wb.Worksheets.Item("PeriodPredictions")
$periodPredictions.Activate()
$periodPredictions.Range("I8").Select()
